$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.149.62'
$ws.Range("E2").Value = '  -0.29%  '

# Row 3
$ws.Range("D3").Value = '1.881.77'
$ws.Range("E3").Value = '  -1.38%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '''313.60'
$ws.Range("E5").Value = '  -0.28%  '

# Row 6
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.10%  '

# Row 7
$ws.Range("D7").Value = '''0.5087'
$ws.Range("E7").Value = '  +0.45%  '

# Row 9
$ws.Range("D9").Value = '''0.09149'
$ws.Range("E9").Value = '  -2.11%  '

# Row 10
$ws.Range("E10").Value = '  -1.51%  '

# Row 11
$ws.Range("D11").Value = '''41.58'
$ws.Range("E11").Value = '  -0.60%  '

# Row 12
$ws.Range("D12").Value = '''6.360'
$ws.Range("E12").Value = '  -0.67%  '

# Row 13
$ws.Range("D13").Value = '''20.81'
$ws.Range("E13").Value = '  -0.35%  '

# Row 14
$ws.Range("D14").Value = '1.875.51'
$ws.Range("E14").Value = '  -2.08%  '

# Row 15
$ws.Range("D15").Value = '''7.212'
$ws.Range("E15").Value = '  -1.31%  '

# Row 16
$ws.Range("E16").Value = '  +0.09%  '

# Row 17
$ws.Range("D17").Value = '''0.00001115'
$ws.Range("E17").Value = '  -0.87%  '

# Row 18
$ws.Range("D18").Value = '''91.28'
$ws.Range("E18").Value = '  -1.55%  '

# Row 19
$ws.Range("D19").Value = '''0.06613'
$ws.Range("E19").Value = '  +0.14%  '

# Row 20
$ws.Range("E20").Value = '  +1.43%  '

# Row 21
$ws.Range("E21").Value = '  +0.13%  '

# Row 22
$ws.Range("D22").Value = '''6.120'
$ws.Range("E22").Value = '  -1.30%  '

# Row 23
$ws.Range("D23").Value = '28.188.21'

# Row 24
$ws.Range("D24").Value = '''11.46'
$ws.Range("E24").Value = '  +0.32%  '

# Row 25
$ws.Range("D25").Value = '''2.278'
$ws.Range("E25").Value = '  -1.75%  '

# Row 26
$ws.Range("D26").Value = '''2.578'
$ws.Range("E26").Value = '  -0.53%  '

# Row 27
$ws.Range("D27").Value = '2.089.59'
$ws.Range("E27").Value = '  -1.94%  '

# Row 28
$ws.Range("D28").Value = '''20.83'
$ws.Range("E28").Value = '  -1.11%  '

# Row 29
$ws.Range("D29").Value = '''157.60'
$ws.Range("E29").Value = '  -0.19%  '

# Row 30
$ws.Range("D30").Value = '''126.90'
$ws.Range("E30").Value = '  -0.29%  '

# Row 31
$ws.Range("D31").Value = '''1.067'
$ws.Range("E31").Value = '  -3.25%  '

# Row 32
$ws.Range("D32").Value = '''0.1055'
$ws.Range("E32").Value = '  -1.66%  '

# Row 33
$ws.Range("D33").Value = '''5.625'
$ws.Range("E33").Value = '  -0.24%  '

# Row 34
$ws.Range("D34").Value = '''3.605'
$ws.Range("E34").Value = '  -0.23%  '

# Row 35
$ws.Range("D35").Value = '''9.707'
$ws.Range("E35").Value = '  +0.33%  '

# Row 36
$ws.Range("D36").Value = '''0.02460'
$ws.Range("E36").Value = '  +1.85%  '

# Row 37
$ws.Range("D37").Value = '''0.06588'
$ws.Range("E37").Value = '  -1.23%  '

# Row 38
$ws.Range("D38").Value = '''0.2180'
$ws.Range("E38").Value = '  -0.35%  '

# Row 39
$ws.Range("D39").Value = '''1.214'
$ws.Range("E39").Value = '  -2.80%  '

# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''1.243'
$ws.Range("E40").Value = '  -2.93%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.6426'
$ws.Range("E41").Value = '  -0.02%  '

# Row 42
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '''11.58'
$ws.Range("E42").Value = '  +0.82%  '

# Row 43
$ws.Range("D43").Value = '''4.938'
$ws.Range("E43").Value = '  -1.51%  '

# Row 44
$ws.Range("D44").Value = '''13.25'
$ws.Range("E44").Value = '  -0.32%  '

# Row 45
$ws.Range("D45").Value = '''0.6020'
$ws.Range("E45").Value = '  +0.18%  '

# Row 46
$ws.Range("E46").Value = '  -1.20%  '

# Row 47
$ws.Range("E47").Value = '  -0.16%  '

# Row 48
$ws.Range("D48").Value = '''1.234'
$ws.Range("E48").Value = '  +4.02%  '

# Row 49
$ws.Range("D49").Value = '''2.005'
$ws.Range("E49").Value = '  -0.73%  '

# Row 50
$ws.Range("D50").Value = '''121.44'
$ws.Range("E50").Value = '  -1.25%  '

# Row 51
$ws.Range("D51").Value = '''80.00'
$ws.Range("E51").Value = '  +2.06%  '
